$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.661.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.99%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4826'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2872'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06557'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.964.29'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07472'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.98%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.108'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.22'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6688'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.641.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("E17").Value = '  +1.08%  '

$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.196.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.59%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007581'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '233.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.281'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.234'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.61%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '169.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.353'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.965'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1023'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.395'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.330'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.036'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05066'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.217'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7551'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.002'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.76%  '

$ws.Range("E38").Value = '  +2.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.648'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9196'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.074'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4297'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '

$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.676'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.450'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1275'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.499'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.003'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.80%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.03%  '
